$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.573.34'
$ws.Range("E2").Value = '  +0.33%  '

$ws.Range("D3").Value = '1.874.40'
$ws.Range("E3").Value = '  -1.01%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.71'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4789'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.23%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2813'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06481'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.04%  '

$ws.Range("D10").Value = '1.944.50'
$ws.Range("E10").Value = '  +2.71%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07470'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.73%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.44'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.50%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.070'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.61%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '87.74'
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6607'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.23%  '

$ws.Range("D16").Value = '30.540.89'
$ws.Range("E16").Value = '  +0.27%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.21'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.54%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.02%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007542'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.21%  '

$ws.Range("B20").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D20").Value = '2.137.16'
$ws.Range("E20").Value = '  -0.09%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '226.86'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.17%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.003'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.12%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.258'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.30%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.127'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.90%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.277'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.32%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.36'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.46'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.83%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.926'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.22%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.400'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.49%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09658'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.33%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.322'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.55%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.994'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05049'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.67%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.209'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7440'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.93%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.709'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.47%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01854'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.07%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.635'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.46%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9088'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.13%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.058'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.56%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '105.67'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.96%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4251'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.52%  '

$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9990'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.45%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.753'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.50%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.310'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.25%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1285'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.48%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '63.81'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.89%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.920'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.12%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.458'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.74%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.60'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.91%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05648'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.40%  '
